$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44466
$ws.Range("M2").Value = 60

# Row 3
$ws.Range("D3").Value = 44434
$ws.Range("L3").Value = "Primera"
$ws.Range("M3").Value = 20
$ws.Range("N3").Value = 20000
$ws.Range("O3").Value = 20000
$ws.Range("P3").Value = 20000
$ws.Range("S3").Value = 2000

# Row 4
$ws.Range("D4").Value = 44435
$ws.Range("M4").Value = 40
$ws.Range("N4").Value = 20000
$ws.Range("O4").Value = 20000
$ws.Range("P4").Value = 20000
$ws.Range("S4").Value = 2000

# Row 5
$ws.Range("D5").Value = 44476
$ws.Range("M5").Value = 120
$ws.Range("N5").Value = 20000
$ws.Range("O5").Value = 20000
$ws.Range("P5").Value = 20000
$ws.Range("S5").Value = 2000

# Row 6 (D6 unchanged)
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 30000
$ws.Range("O6").Value = 30000
$ws.Range("P6").Value = 30000
$ws.Range("S6").Value = 3000

# Row 7
$ws.Range("D7").Value = 44503
$ws.Range("L7").Value = "Segunda"
$ws.Range("M7").Value = 50
$ws.Range("N7").Value = 25000
$ws.Range("O7").Value = 25000
$ws.Range("P7").Value = 25000
$ws.Range("S7").Value = 2500

# Row 8
$ws.Range("D8").Value = 44517
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 100
$ws.Range("N8").Value = 27000
$ws.Range("O8").Value = 27000
$ws.Range("P8").Value = 27000
$ws.Range("S8").Value = 2700

# Row 9
$ws.Range("D9").Value = 44517
$ws.Range("M9").Value = 30
$ws.Range("N9").Value = 25000
$ws.Range("O9").Value = 25000
$ws.Range("P9").Value = 25000
$ws.Range("S9").Value = 2500

# Row 10
$ws.Range("D10").Value = 44473
$ws.Range("M10").Value = 180

# Row 12
$ws.Range("D12").Value = 44432
$ws.Range("M12").Value = 20
